# ----------------------------------------------------------------------
# Applies the diff:
#   1. Inserts 3 new paragraphs at the very top of the document body:
#        - "Scientific Background" (bold heading, sz24)
#        - "Detecting deception ... used as well." (sz24, not bold, 5 runs)
#        - an empty paragraph (sz24)
#   2. Rewrites the "Some approaches ..." paragraph later in the doc,
#      splitting it into several runs with a handful of wording tweaks.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-RunsByBookmark($rng, [int[]]$offsets) {
    # Forces run boundaries inside $rng (a Range already containing the
    # final text) at the given character offsets (relative to $rng.Start)
    # by briefly wrapping a Bookmark around the leading chunk and then
    # deleting the bookmark again -- this leaves the run split in place
    # without leaving any bookmark markup behind.
    $i = 0
    foreach ($off in $offsets) {
        $i = $i + 1
        $bmName = "tmpSplitBm" + $i
        $sub = $d.Range($rng.Start, $rng.Start + $off)
        $d.Bookmarks.Add($bmName, $sub) | Out-Null
        $d.Bookmarks($bmName).Delete()
    }
}

# ------------------------------------------------------------------
# Part 1: three new paragraphs at the top of the document.
# ------------------------------------------------------------------

# --- Paragraph 1: bold heading "Scientific Background" ---
# Copy paragraph 1 ("Preprocessing") so we inherit its exact bold/sz24
# paragraph-mark + run formatting, paste a duplicate before it, then
# overwrite its text.
$headingSrc = $d.Paragraphs(1)
$headingSrc.Range.Copy()
$insPoint = $d.Range(0, 0)
$insPoint.Paste()
$d.Paragraphs(1).Range.Text = "Scientific Background"

# --- Paragraph 2: body paragraph, sz24, not bold, 5 runs ---
# Build it at the very end of the document (a location with no bold
# formatting in context), then cut it and paste it into place so it
# does not inherit the heading's bold run properties.
$bodyText = "Detecting deception through audio involves analyzing subtle vocal cues such as pitch, tone, and speech pace. Mel-Frequency Cepstral Coefficients (MFCCs) are key in quantifying these spectral properties of speech, highlighting subtle fluctuations in vocal expression. Bidirectional Long Short-Term Memory (Bi-LSTM) networks, a type of Recurrent Neural Network (RNN), excel in processing these features by examining speech patterns from both past and future contexts. This method allows for the detection of inconsistencies and anomalies in speech that typically indicate lying. RNNs in general are well suited for this task, although simple machine learning algorithms like SVMs or Random forests can be used as well."

$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newBody = $d.Paragraphs($n).Range
$newBody.Font.Size = 12
$newBody.Font.SizeBi = 12
$newBody.Text = $bodyText

# Split into the 5 runs shown by the diff:
#   "...indicate lying" | "." | " RNNs in general" | " are well suited for this task" | ", although ... used as well."
$bodyRng = $d.Paragraphs($n).Range
Split-RunsByBookmark $bodyRng @(575, 576, 592, 622)

# Move this paragraph to the top, right after "Scientific Background".
$n = $d.Paragraphs.Count
$cutRng = $d.Paragraphs($n).Range
$cutRng.Cut()
$p2Start = $d.Paragraphs(2).Range.Start
$pasteTarget = $d.Range($p2Start, $p2Start)
$pasteTarget.Paste()

# --- Paragraph 3: empty paragraph, sz24 ---
# Same trick: copy paragraph 2 (already sz24 / not bold) so the empty
# paragraph inherits the exact same paragraph-mark formatting, then
# blank out its text.
$fmtSrc = $d.Paragraphs(2)
$fmtSrc.Range.Copy()
$p3Start = $d.Paragraphs(3).Range.Start
$pasteTarget3 = $d.Range($p3Start, $p3Start)
$pasteTarget3.Paste()
$p3 = $d.Paragraphs(3).Range
$textOnly = $d.Range($p3.Start, $p3.End - 1)
$textOnly.Delete()

# ------------------------------------------------------------------
# Part 2: rewrite the "Some approaches ..." paragraph.
# ------------------------------------------------------------------

$oldSentence = "Some approaches promote splitting every single audio signal into several ones and compute MFCC of each, but these approaches are not accurate because it is not clear which part of the clip displayed deception cues. (although this approach does act as some sort of augmentation method, but its results are illogical and inaccurate)"
$newSentence = "Some approaches suggest splitting every single audio signal into several smaller ones and computing the MFCCs of each, but these approaches are not accurate because it is not clear which part of the clip displayed deception cues. (although this approach does achieve some sort of augmentation especially considering the limited amount of data available, but its results are illogical and inaccurate)"

$found = $d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# Find the paragraph that now holds the replaced sentence so we can
# split it into the runs shown by the diff.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Some approaches suggest")) {
        $target = $para
        break
    }
}

$tRng = $target.Range
Split-RunsByBookmark $tRng @(16, 23, 73, 81, 96, 103, 108, 109, 259, 266, 293, 303, 352)

Write-Output "done"
